# Adds two new columns, I ("I0") and J ("IF"), to the stats sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Give the new header cells (I1, J1) the same formatting (bold font,
# border, centered/top alignment) as the existing header cells by
# copying the format from H1.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# Header labels
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Row 2 holds special values that don't follow the rest of the table.
$ws.Range("I2").Value = 7
$ws.Range("J2").Value = 9

# For every remaining data row, I is always 1 and J mirrors the
# existing IP value already present in column H.
for ($r = 3; $r -le 35; $r++) {
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $ws.Cells.Item($r, 8).Value2
}
